$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as text so values like
# "29.220.13" or "1.000" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.220.13"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3
$ws.Range("D3").Value = "1.840.00"
$ws.Range("E3").Value = "  -0.03%  "

# Row 4
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "241.07"
$ws.Range("E5").Value = "  -0.96%  "

# Row 6
$ws.Range("D6").Value = "0.6702"
$ws.Range("E6").Value = "  -2.50%  "

# Row 8
$ws.Range("D8").Value = "0.07416"
$ws.Range("E8").Value = "  -1.54%  "

# Row 9
$ws.Range("D9").Value = "0.2962"
$ws.Range("E9").Value = "  -2.25%  "

# Row 10
$ws.Range("D10").Value = "22.75"
$ws.Range("E10").Value = "  -2.09%  "

# Row 11
$ws.Range("D11").Value = "0.07720"
$ws.Range("E11").Value = "  +0.58%  "

# Row 12
$ws.Range("D12").Value = "5.022"
$ws.Range("E12").Value = "  -1.34%  "

# Row 13
$ws.Range("D13").Value = "0.6773"
$ws.Range("E13").Value = "  -1.46%  "

# Row 14
$ws.Range("D14").Value = "1.735.15"
$ws.Range("E14").Value = "  -5.71%  "

# Row 15
$ws.Range("D15").Value = "86.39"
$ws.Range("E15").Value = "  -2.88%  "

# Row 16
$ws.Range("D16").Value = "6.167"
$ws.Range("E16").Value = "  -1.88%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.000008271"
$ws.Range("E17").Value = "  +0.43%  "

# Row 18
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "28.835.77"
$ws.Range("E18").Value = "  -1.62%  "

# Row 19
$ws.Range("D19").Value = "228.37"
$ws.Range("E19").Value = "  -2.00%  "

# Row 20
$ws.Range("E20").Value = "  -0.49%  "

# Row 21
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("D22").Value = "7.154"
$ws.Range("E22").Value = "  -4.63%  "

# Row 23
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("D24").Value = "159.89"
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("D25").Value = "8.694"
$ws.Range("E25").Value = "  -1.69%  "

# Row 26
$ws.Range("D26").Value = "0.1403"
$ws.Range("E26").Value = "  -4.02%  "

# Row 27
$ws.Range("D27").Value = "18.02"
$ws.Range("E27").Value = "  -0.67%  "

# Row 28
$ws.Range("D28").Value = "1.511"
$ws.Range("E28").Value = "  -0.39%  "

# Row 29
$ws.Range("D29").Value = "4.196"
$ws.Range("E29").Value = "  -1.06%  "

# Row 30
$ws.Range("D30").Value = "4.084"
$ws.Range("E30").Value = "  -1.55%  "

# Row 31
$ws.Range("D31").Value = "1.193"
$ws.Range("E31").Value = "  -0.98%  "

# Row 32
$ws.Range("D32").Value = "0.05364"
$ws.Range("E32").Value = "  +3.76%  "

# Row 33
$ws.Range("D33").Value = "1.875"
$ws.Range("E33").Value = "  +1.73%  "

# Row 34
$ws.Range("D34").Value = "0.7567"
$ws.Range("E34").Value = "  -1.91%  "

# Row 35
$ws.Range("D35").Value = "1.141"
$ws.Range("E35").Value = "  +0.20%  "

# Row 36
$ws.Range("D36").Value = "2.677"

# Row 37
$ws.Range("D37").Value = "1.328.98"
$ws.Range("E37").Value = "  +2.43%  "

# Row 38
$ws.Range("D38").Value = "0.01803"
$ws.Range("E38").Value = "  -2.36%  "

# Row 39
$ws.Range("D39").Value = "2.728"
$ws.Range("E39").Value = "  +0.93%  "

# Row 40
$ws.Range("D40").Value = "0.9224"
$ws.Range("E40").Value = "  -2.68%  "

# Row 41
$ws.Range("D41").Value = "5.953"
$ws.Range("E41").Value = "  +2.37%  "

# Row 42
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.26%  "

# Row 43
$ws.Range("D43").Value = "103.11"
$ws.Range("E43").Value = "  -1.91%  "

# Row 44
$ws.Range("D44").Value = "0.08197"
$ws.Range("E44").Value = "  +16.56%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.00000000124"
$ws.Range("E45").Value = "  +0.99%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.5161"
$ws.Range("E46").Value = "  -0.74%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "1.765"
$ws.Range("E47").Value = "  -0.60%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "63.65"
$ws.Range("E48").Value = "  -1.96%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.284"
$ws.Range("E49").Value = "  -4.11%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "1.870.82"
$ws.Range("E50").Value = "  -5.72%  "

# Row 51
$ws.Range("D51").Value = "0.05921"
$ws.Range("E51").Value = "  -0.07%  "
